$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/boiling-flasks-round-bottom/', 'http://drupal-geneseo-backup.ddev.site/themes/custom/geneseo_bootstrap5/images/geneseo-logo-simple.png', 'image', '/themes/custom/geneseo_bootstrap5/images/geneseo-logo-simple.png', 200),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/boiling-flasks-round-bottom/', 'http://drupal-geneseo-backup.ddev.site/themes/custom/geneseo_bootstrap5/images/my-geneseo.png', 'image', '/themes/custom/geneseo_bootstrap5/images/my-geneseo.png', 200),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/boiling-flasks-round-bottom/', 'http://drupal-geneseo-backup.ddev.site/sites/default/files/sites/404%20Page/page-not-found-graphic.jpg', 'image', '/sites/default/files/sites/404%20Page/page-not-found-graphic.jpg', 404),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/boiling-flasks-round-bottom/', 'https://www.geneseo.edu/sites/default/files/images/geneseo-logo.svg', 'image', 'https://www.geneseo.edu/sites/default/files/images/geneseo-logo.svg', 200),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/boiling-flasks-round-bottom/', 'https://www.geneseo.edu/sites/default/files/sites/social/facebook.png', 'image', 'https://www.geneseo.edu/sites/default/files/sites/social/facebook.png', 200),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/boiling-flasks-round-bottom/', 'https://www.geneseo.edu/sites/default/files/sites/social/instagram.png', 'image', 'https://www.geneseo.edu/sites/default/files/sites/social/instagram.png', 200),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/boiling-flasks-round-bottom/', 'https://www.geneseo.edu/sites/default/files/sites/social/linkedin.png', 'image', 'https://www.geneseo.edu/sites/default/files/sites/social/linkedin.png', 200),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/boiling-flasks-round-bottom/', 'https://www.geneseo.edu/sites/default/files/images/bluesky.png', 'image', 'https://www.geneseo.edu/sites/default/files/images/bluesky.png', 200),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/boiling-flasks-round-bottom/', 'https://www.geneseo.edu/sites/default/files/sites/social/youtube.png', 'image', 'https://www.geneseo.edu/sites/default/files/sites/social/youtube.png', 200),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/boiling-flasks-round-bottom/', 'https://bat.bing.com/action/0?ti=73008797&tm=gtm002&Ver=2&mid=3aa57b59-4993-466e-bd35-d7d65d5121d7&bo=1&sid=858ae0709e0811f092bf639d456feb70&vid=858b4a509e0811f0973b254feeb8ff2e&vids=1&msclkid=N&gtm_tag_source=1&pi=0&lg=en-US&sw=1280&sh=720&sc=24&nwd=1&tl=Page%20Not%20Found%20%7C%20SUNY%20Geneseo&p=http%3A%2F%2Fdrupal-geneseo-backup.ddev.site%2Fbasic_page%2Fboiling-flasks-round-bottom%2F&r=&lt=489&evt=pageLoad&sv=2&cdb=AQAQ&rn=592926', 'image', 'https://bat.bing.com/action/0?ti=73008797&tm=gtm002&Ver=2&mid=3aa57b59-4993-466e-bd35-d7d65d5121d7&bo=1&sid=858ae0709e0811f092bf639d456feb70&vid=858b4a509e0811f0973b254feeb8ff2e&vids=1&msclkid=N&gtm_tag_source=1&pi=0&lg=en-US&sw=1280&sh=720&sc=24&nwd=1&tl=Page%20Not%20Found%20%7C%20SUNY%20Geneseo&p=http%3A%2F%2Fdrupal-geneseo-backup.ddev.site%2Fbasic_page%2Fboiling-flasks-round-bottom%2F&r=&lt=489&evt=pageLoad&sv=2&cdb=AQAQ&rn=592926', 204),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/journals-and-reflections/', 'http://drupal-geneseo-backup.ddev.site/themes/custom/geneseo_bootstrap5/images/geneseo-logo-simple.png', 'image', '/themes/custom/geneseo_bootstrap5/images/geneseo-logo-simple.png', 200),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/journals-and-reflections/', 'http://drupal-geneseo-backup.ddev.site/themes/custom/geneseo_bootstrap5/images/my-geneseo.png', 'image', '/themes/custom/geneseo_bootstrap5/images/my-geneseo.png', 200),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/journals-and-reflections/', 'http://drupal-geneseo-backup.ddev.site/sites/default/files/sites/404%20Page/page-not-found-graphic.jpg', 'image', '/sites/default/files/sites/404%20Page/page-not-found-graphic.jpg', 404),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/journals-and-reflections/', 'https://www.geneseo.edu/sites/default/files/images/geneseo-logo.svg', 'image', 'https://www.geneseo.edu/sites/default/files/images/geneseo-logo.svg', 200),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/journals-and-reflections/', 'https://www.geneseo.edu/sites/default/files/sites/social/facebook.png', 'image', 'https://www.geneseo.edu/sites/default/files/sites/social/facebook.png', 200),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/journals-and-reflections/', 'https://www.geneseo.edu/sites/default/files/sites/social/instagram.png', 'image', 'https://www.geneseo.edu/sites/default/files/sites/social/instagram.png', 200),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/journals-and-reflections/', 'https://www.geneseo.edu/sites/default/files/sites/social/linkedin.png', 'image', 'https://www.geneseo.edu/sites/default/files/sites/social/linkedin.png', 200),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/journals-and-reflections/', 'https://www.geneseo.edu/sites/default/files/images/bluesky.png', 'image', 'https://www.geneseo.edu/sites/default/files/images/bluesky.png', 200),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/journals-and-reflections/', 'https://www.geneseo.edu/sites/default/files/sites/social/youtube.png', 'image', 'https://www.geneseo.edu/sites/default/files/sites/social/youtube.png', 200),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/journals-and-reflections/', 'https://bat.bing.com/action/0?ti=73008797&tm=gtm002&Ver=2&mid=6d0e4dec-ba53-4693-ae40-b56f68c60c12&bo=1&sid=859cb4d09e0811f0aa011f88ece24a43&vid=859ce0609e0811f096f3ad0e0bb56ae6&vids=1&msclkid=N&gtm_tag_source=1&pi=0&lg=en-US&sw=1280&sh=720&sc=24&nwd=1&tl=Page%20Not%20Found%20%7C%20SUNY%20Geneseo&p=http%3A%2F%2Fdrupal-geneseo-backup.ddev.site%2Fbasic_page%2Fjournals-and-reflections%2F&r=&lt=604&evt=pageLoad&sv=2&cdb=AQAQ&rn=837458', 'image', 'https://bat.bing.com/action/0?ti=73008797&tm=gtm002&Ver=2&mid=6d0e4dec-ba53-4693-ae40-b56f68c60c12&bo=1&sid=859cb4d09e0811f0aa011f88ece24a43&vid=859ce0609e0811f096f3ad0e0bb56ae6&vids=1&msclkid=N&gtm_tag_source=1&pi=0&lg=en-US&sw=1280&sh=720&sc=24&nwd=1&tl=Page%20Not%20Found%20%7C%20SUNY%20Geneseo&p=http%3A%2F%2Fdrupal-geneseo-backup.ddev.site%2Fbasic_page%2Fjournals-and-reflections%2F&r=&lt=604&evt=pageLoad&sv=2&cdb=AQAQ&rn=837458', 204),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/conh35clcl2-prep/', 'http://drupal-geneseo-backup.ddev.site/themes/custom/geneseo_bootstrap5/images/geneseo-logo-simple.png', 'image', '/themes/custom/geneseo_bootstrap5/images/geneseo-logo-simple.png', 200),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/conh35clcl2-prep/', 'http://drupal-geneseo-backup.ddev.site/themes/custom/geneseo_bootstrap5/images/my-geneseo.png', 'image', '/themes/custom/geneseo_bootstrap5/images/my-geneseo.png', 200),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/conh35clcl2-prep/', 'http://drupal-geneseo-backup.ddev.site/sites/default/files/sites/404%20Page/page-not-found-graphic.jpg', 'image', '/sites/default/files/sites/404%20Page/page-not-found-graphic.jpg', 404),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/conh35clcl2-prep/', 'https://www.geneseo.edu/sites/default/files/images/geneseo-logo.svg', 'image', 'https://www.geneseo.edu/sites/default/files/images/geneseo-logo.svg', 200),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/conh35clcl2-prep/', 'https://www.geneseo.edu/sites/default/files/sites/social/facebook.png', 'image', 'https://www.geneseo.edu/sites/default/files/sites/social/facebook.png', 200),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/conh35clcl2-prep/', 'https://www.geneseo.edu/sites/default/files/sites/social/instagram.png', 'image', 'https://www.geneseo.edu/sites/default/files/sites/social/instagram.png', 200),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/conh35clcl2-prep/', 'https://www.geneseo.edu/sites/default/files/sites/social/linkedin.png', 'image', 'https://www.geneseo.edu/sites/default/files/sites/social/linkedin.png', 200),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/conh35clcl2-prep/', 'https://www.geneseo.edu/sites/default/files/images/bluesky.png', 'image', 'https://www.geneseo.edu/sites/default/files/images/bluesky.png', 200),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/conh35clcl2-prep/', 'https://www.geneseo.edu/sites/default/files/sites/social/youtube.png', 'image', 'https://www.geneseo.edu/sites/default/files/sites/social/youtube.png', 200),
    @('2-2', 'http://drupal-geneseo-backup.ddev.site/basic_page/conh35clcl2-prep/', 'https://bat.bing.com/action/0?ti=73008797&tm=gtm002&Ver=2&mid=efdde95e-307b-4b01-8286-6279939ba5e9&bo=1&sid=859cb4d09e0811f0aa011f88ece24a43&vid=859ce0609e0811f096f3ad0e0bb56ae6&vids=0&msclkid=N&gtm_tag_source=1&pi=0&lg=en-US&sw=1280&sh=720&sc=24&nwd=1&tl=Page%20Not%20Found%20%7C%20SUNY%20Geneseo&p=http%3A%2F%2Fdrupal-geneseo-backup.ddev.site%2Fbasic_page%2Fconh35clcl2-prep%2F&r=&lt=255&evt=pageLoad&sv=2&cdb=AQkR&rn=828266', 'image', 'https://bat.bing.com/action/0?ti=73008797&tm=gtm002&Ver=2&mid=efdde95e-307b-4b01-8286-6279939ba5e9&bo=1&sid=859cb4d09e0811f0aa011f88ece24a43&vid=859ce0609e0811f096f3ad0e0bb56ae6&vids=0&msclkid=N&gtm_tag_source=1&pi=0&lg=en-US&sw=1280&sh=720&sc=24&nwd=1&tl=Page%20Not%20Found%20%7C%20SUNY%20Geneseo&p=http%3A%2F%2Fdrupal-geneseo-backup.ddev.site%2Fbasic_page%2Fconh35clcl2-prep%2F&r=&lt=255&evt=pageLoad&sv=2&cdb=AQkR&rn=828266', 204),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/mofs-part-3/', 'http://drupal-geneseo-backup.ddev.site/themes/custom/geneseo_bootstrap5/images/geneseo-logo-simple.png', 'image', '/themes/custom/geneseo_bootstrap5/images/geneseo-logo-simple.png', 200),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/mofs-part-3/', 'http://drupal-geneseo-backup.ddev.site/themes/custom/geneseo_bootstrap5/images/my-geneseo.png', 'image', '/themes/custom/geneseo_bootstrap5/images/my-geneseo.png', 200),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/mofs-part-3/', 'http://drupal-geneseo-backup.ddev.site/sites/default/files/sites/404%20Page/page-not-found-graphic.jpg', 'image', '/sites/default/files/sites/404%20Page/page-not-found-graphic.jpg', 404),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/mofs-part-3/', 'https://www.geneseo.edu/sites/default/files/images/geneseo-logo.svg', 'image', 'https://www.geneseo.edu/sites/default/files/images/geneseo-logo.svg', 200),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/mofs-part-3/', 'https://www.geneseo.edu/sites/default/files/sites/social/facebook.png', 'image', 'https://www.geneseo.edu/sites/default/files/sites/social/facebook.png', 200),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/mofs-part-3/', 'https://www.geneseo.edu/sites/default/files/sites/social/instagram.png', 'image', 'https://www.geneseo.edu/sites/default/files/sites/social/instagram.png', 200),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/mofs-part-3/', 'https://www.geneseo.edu/sites/default/files/sites/social/linkedin.png', 'image', 'https://www.geneseo.edu/sites/default/files/sites/social/linkedin.png', 200),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/mofs-part-3/', 'https://www.geneseo.edu/sites/default/files/images/bluesky.png', 'image', 'https://www.geneseo.edu/sites/default/files/images/bluesky.png', 200),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/mofs-part-3/', 'https://www.geneseo.edu/sites/default/files/sites/social/youtube.png', 'image', 'https://www.geneseo.edu/sites/default/files/sites/social/youtube.png', 200),
    @('2-1', 'http://drupal-geneseo-backup.ddev.site/basic_page/mofs-part-3/', 'https://bat.bing.com/action/0?ti=73008797&tm=gtm002&Ver=2&mid=0af70e59-1947-4c56-81b2-8eb0399133d7&bo=1&sid=858ae0709e0811f092bf639d456feb70&vid=858b4a509e0811f0973b254feeb8ff2e&vids=0&msclkid=N&gtm_tag_source=1&pi=0&lg=en-US&sw=1280&sh=720&sc=24&nwd=1&tl=Page%20Not%20Found%20%7C%20SUNY%20Geneseo&p=http%3A%2F%2Fdrupal-geneseo-backup.ddev.site%2Fbasic_page%2Fmofs-part-3%2F&r=&lt=180&evt=pageLoad&sv=2&cdb=AQkR&rn=970059', 'image', 'https://bat.bing.com/action/0?ti=73008797&tm=gtm002&Ver=2&mid=0af70e59-1947-4c56-81b2-8eb0399133d7&bo=1&sid=858ae0709e0811f092bf639d456feb70&vid=858b4a509e0811f0973b254feeb8ff2e&vids=0&msclkid=N&gtm_tag_source=1&pi=0&lg=en-US&sw=1280&sh=720&sc=24&nwd=1&tl=Page%20Not%20Found%20%7C%20SUNY%20Geneseo&p=http%3A%2F%2Fdrupal-geneseo-backup.ddev.site%2Fbasic_page%2Fmofs-part-3%2F&r=&lt=180&evt=pageLoad&sv=2&cdb=AQkR&rn=970059', 204)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r = $r + 1
}

$ws.Rows("42:60").Delete()